$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find a paragraph by a unique snippet of its text and replace the
# whole paragraph's XML. Used both to strip the green "w:color 008000"
# direct formatting that used to mark some rows as translatable notes, and
# to edit text while leaving every other attribute (rsids, paraId, ...)
# untouched.
# ---------------------------------------------------------------------------
function Set-ParagraphXml([string]$findText, [string]$newParaXml) {
    $r = $d.Content
    $r.Find.ClearFormatting()
    $ok = $r.Find.Execute($findText)
    if (-not $ok) {
        Write-Host "NOT FOUND: $findText"
        return
    }
    $p = $r.Paragraphs(1)
    $pr = $p.Range
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>' + $newParaXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $pr.InsertXML($pkg) | Out-Null
}

# ---------------------------------------------------------------------------
# 1) "Please make at least one connection to ground (triangle symbol)"
#    -> "Please make at least one connection to ground"
#    (plain text edit, no formatting change)
# ---------------------------------------------------------------------------
$p0 = '<w:p w14:paraId="5AF5C7B7" w14:textId="77777777" w:rsidR="003C3C64" w:rsidRPr="00BE43D7" w:rsidRDefault="003C3C64" w:rsidP="003C3C64"><w:pPr><w:pStyle w:val="PlainText"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r w:rsidRPr="00BE43D7"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>Please make at least one connection to ground</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:hint="eastAsia"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>.</w:t></w:r></w:p>'
Set-ParagraphXml "Please make at least one connection to ground (triangle symbol)" $p0

# ---------------------------------------------------------------------------
# 2) "Warning! Simulator might produce meaningless results or no result
#    with illegal circuits." - drop the green color from the paragraph
#    mark and both runs.
# ---------------------------------------------------------------------------
$p1 = '<w:p w14:paraId="2C8EE634" w14:textId="16B73A85" w:rsidR="003C3C64" w:rsidRPr="003B06FD" w:rsidRDefault="00FC32E9" w:rsidP="00607FB4"><w:pPr><w:pStyle w:val="PlainText"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr><w:r w:rsidRPr="003B06FD"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>Warning!</w:t></w:r><w:r w:rsidR="003C3C64" w:rsidRPr="003B06FD"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> Simulator might produce meaningless results or no result with illegal circuits.</w:t></w:r></w:p>'
Set-ParagraphXml "Warning! Simulator might produce meaningless results or no result with illegal circuits." $p1

# ---------------------------------------------------------------------------
# 3) "Newton Method failed; do your current sources have a conductive path
#    to ground?" - drop the green color from the paragraph mark and all
#    three runs.
# ---------------------------------------------------------------------------
$p2 = '<w:p w14:paraId="5148C86C" w14:textId="40739FD6" w:rsidR="003C3C64" w:rsidRPr="003B06FD" w:rsidRDefault="003C3C64" w:rsidP="00607FB4"><w:pPr><w:pStyle w:val="PlainText"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr><w:r w:rsidRPr="003B06FD"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">Newton Method </w:t></w:r><w:r w:rsidR="003B06FD" w:rsidRPr="003B06FD"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>failed;</w:t></w:r><w:r w:rsidRPr="003B06FD"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> do your current sources have a conductive path to ground?</w:t></w:r></w:p>'
Set-ParagraphXml "Newton Method failed; do your current sources have a conductive path to ground?" $p2

# ---------------------------------------------------------------------------
# 4) "Newton Method failed; it may be your circuit or it may be our
#    simulator." - drop the green color from the paragraph mark and both
#    runs.
# ---------------------------------------------------------------------------
$p3 = '<w:p w14:paraId="70811FBD" w14:textId="3A4B40CA" w:rsidR="003C3C64" w:rsidRPr="003B06FD" w:rsidRDefault="003B06FD" w:rsidP="00607FB4"><w:pPr><w:pStyle w:val="PlainText"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>Newton Method failed;</w:t></w:r><w:r w:rsidR="003C3C64" w:rsidRPr="003B06FD"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> it may be your circuit or it may be our simulator.</w:t></w:r></w:p>'
Set-ParagraphXml "Newton Method failed; it may be your circuit or it may be our simulator." $p3

Write-Host "Done"
